$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.908.88"
$ws.Range("E2").Value = "'  +5.33%  "

$ws.Range("D3").Value = "'3.352.91"
$ws.Range("E3").Value = "'  +5.23%  "

$ws.Range("E4").Value = "'  +0.04%  "

$ws.Range("D5").Value = "'570.77"
$ws.Range("E5").Value = "'  +6.76%  "

$ws.Range("D6").Value = "'152.49"
$ws.Range("E6").Value = "'  +5.52%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.00%  "

$ws.Range("D8").Value = "'3.357.04"
$ws.Range("E8").Value = "'  +5.20%  "

$ws.Range("E9").Value = "'  -0.34%  "

$ws.Range("D10").Value = "'7.43"
$ws.Range("E10").Value = "'  +1.56%  "

$ws.Range("E11").Value = "'  +4.97%  "

$ws.Range("D12").Value = "'0.438"
$ws.Range("E12").Value = "'  +2.25%  "

$ws.Range("D13").Value = "'3.932.11"
$ws.Range("E13").Value = "'  +5.23%  "

$ws.Range("E14").Value = "'  -0.14%  "

$ws.Range("B15").Value = "'Avalanche"
$ws.Range("C15").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'26.94"
$ws.Range("E15").Value = "'  +3.75%  "

$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000180"
$ws.Range("E16").Value = "'  +4.45%  "

$ws.Range("D17").Value = "'62.938.09"
$ws.Range("E17").Value = "'  +5.30%  "

$ws.Range("D18").Value = "'3.353.83"
$ws.Range("E18").Value = "'  +6.37%  "

$ws.Range("D19").Value = "'6.33"
$ws.Range("E19").Value = "'  +2.03%  "

$ws.Range("D20").Value = "'13.84"
$ws.Range("E20").Value = "'  +5.37%  "

$ws.Range("E21").Value = "'  +2.89%  "

$ws.Range("D22").Value = "'384.61"
$ws.Range("E22").Value = "'  +4.95%  "

$ws.Range("E23").Value = "'  +0.17%  "

$ws.Range("D24").Value = "'0.532"
$ws.Range("E24").Value = "'  +2.56%  "

$ws.Range("D25").Value = "'70.42"
$ws.Range("E25").Value = "'  +1.30%  "

$ws.Range("D26").Value = "'9.45"
$ws.Range("E26").Value = "'  +7.75%  "

$ws.Range("E27").Value = "'  +6.38%  "

$ws.Range("D28").Value = "'0.0₃0961"
$ws.Range("E28").Value = "'  +9.32%  "

$ws.Range("E29").Value = "'  -0.05%  "

$ws.Range("E30").Value = "'  +6.45%  "

$ws.Range("D31").Value = "'22.96"
$ws.Range("E31").Value = "'  +3.27%  "

$ws.Range("D32").Value = "'5.57"
$ws.Range("E32").Value = "'  +5.57%  "

$ws.Range("B33").Value = "'Fetch.AI"
$ws.Range("C33").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.31"
$ws.Range("E33").Value = "'  +9.97%  "

$ws.Range("B34").Value = "'RenderToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'6.28"
$ws.Range("E34").Value = "'  +3.33%  "

$ws.Range("D35").Value = "'6.70"
$ws.Range("E35").Value = "'  +2.31%  "

$ws.Range("E36").Value = "'  +9.77%  "

$ws.Range("D37").Value = "'157.68"
$ws.Range("E37").Value = "'  +0.62%  "

$ws.Range("E38").Value = "'  +12.08%  "

$ws.Range("D39").Value = "'27.01"
$ws.Range("E39").Value = "'  +4.29%  "

$ws.Range("D40").Value = "'0.0331"
$ws.Range("E40").Value = "'  +13.97%  "

$ws.Range("D41").Value = "'0.0737"
$ws.Range("E41").Value = "'  +5.78%  "

$ws.Range("D42").Value = "'2.818.98"
$ws.Range("E42").Value = "'  +1.49%  "

$ws.Range("E43").Value = "'  +3.82%  "

$ws.Range("D44").Value = "'4.27"
$ws.Range("E44").Value = "'  +1.42%  "

$ws.Range("D45").Value = "'0.745"
$ws.Range("E45").Value = "'  +4.60%  "

$ws.Range("D46").Value = "'1.03"
$ws.Range("E46").Value = "'  +5.67%  "

$ws.Range("D47").Value = "'3.397.66"
$ws.Range("E47").Value = "'  +5.31%  "

$ws.Range("D48").Value = "'21.96"
$ws.Range("E48").Value = "'  +7.51%  "

$ws.Range("E49").Value = "'  -2.12%  "

$ws.Range("D50").Value = "'6.31"
$ws.Range("E50").Value = "'  +3.13%  "

$ws.Range("D51").Value = "'292.06"
$ws.Range("E51").Value = "'  +11.17%  "
